$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.100.63"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.790.91"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.546"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.24"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.294"
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0693"
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0941"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "2.048.42"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.50"
$ws.Range("E13").Value = "  +4.02%  "
$ws.Range("D14").Value = "1.781.83"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "34.081.49"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.622"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.20"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.90"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.61"
$ws.Range("D20").Value = "0.0₃0782"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.92"
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.12"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.81"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.18"
$ws.Range("E26").Value = "  +1.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.31"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("E33").Value = "  +3.11%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "1.447.22"
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.646"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  +2.53%  "
$ws.Range("E38").Value = "  +7.33%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.50"
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.928"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.36"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.31"
$ws.Range("E44").Value = "  +6.59%  "
$ws.Range("E45").Value = "  +4.15%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0509"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0138"
$ws.Range("E47").Value = "  -2.89%  "
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.81"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "1.949.99"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  +0.06%  "
